$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was inserted before the existing row 530, which
# pushes every following record (old rows 530-639) down by one row and
# extends the used range to row 640.
$ws.Rows(530).Insert()

# Populate the newly inserted row 530 with the new observation's data.
$ws.Range("A530").Value = 6
$ws.Range("B530").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C530").Value = "Metropolitana"
$ws.Range("D530").Value = 44995
$ws.Range("E530").Value = 13
$ws.Range("F530").Value = 100112039
$ws.Range("G530").Value = "Ciboulette"
$ws.Range("H530").Value = "Sin especificar"
$ws.Range("I530").Value = "Primera"
$ws.Range("J530").Value = 620
$ws.Range("K530").Value = 1500
$ws.Range("L530").Value = 1600
$ws.Range("M530").Value = 1560
$ws.Range("N530").Value = "`$/docena de atados"
$ws.Range("O530").Value = "Región Metropolitana"
$ws.Range("P530").Value = 520
$ws.Range("Q530").Value = 3
$ws.Range("R530").Value = "Hortaliza"
